# Update driver cell values (B6: 50 -> 1000, B7: 2 -> 10)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = 1000
$ws.Range("B7").Value = 10

# Update the rich-text explanatory comment in E1 (merged cell E1:L14)
$ws.Range("E1").Value = "This workbook shows how simple generator functions can stream data into Excel using xlSlim.`nThe Python module simple_counters.py has two counter functions. One is a regular Python function with a yield, which turns the function into a generator function. The second is a slightly more advanced Python function that asynchronously generates values.`nAlways use Python async functions where possible, they are faster.`nCells C6 and D6 show how the functions are exposed in Excel as regular worksheet functions, all the complexity of streaming data is taken care of internally within xlSlim.`nCells A10 and C10 show how numpy arrays can be streamed into Excel. Indeed any supported objects can be streamed. Cells B13 and D13 show how the object cache handles can be passed into xlSlim functions just the same as other xlSlim functions returning cached objects."

# Bold the "Always use Python async functions..." sentence
$boldChars = $ws.Range("E1").Characters(349, 66)
$boldChars.Font.Bold = $true

# Touch the trailing paragraph so it also carries explicit (non-bold) run
# formatting, matching how Excel writes out the run that follows a
# manually-bolded run.
$tailChars = $ws.Range("E1").Characters(415, 440)
$tailChars.Font.Name = "Calibri"

# Row 1 keeps its original (manually set) height even though the comment
# text grew by a line.
$ws.Rows.Item(1).RowHeight = 34.75

# Move the active selection to B6
[void]$ws.Range("B6").Select()
